$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Set the Title value (row 5, column B) to match the Name value (row 4, column B)
$ws.Range("B5").Value = $ws.Range("B4").Text

# Update the Date value (row 8, column B)
$ws.Range("B8").Value = "2024-04-08T12:44:22+00:00"
